{"js": "// Update the worksheet date and every division problem's operands.\n// Each entry is [oldText, newText]; oldText values are unique in the\n// document, so an exact, case-sensitive, non-wildcard search safely\n// identifies the single run to replace.\nconst replacements = [[\"2024-08-04 Sunday\", \"2024-08-05 Monday\"], [\"864\u00f76=\", \"810\u00f77=\"], [\"252\u00f79=\", \"447\u00f78=\"], [\"826\u00f79=\", \"713\u00f77=\"], [\"390\u00f75=\", \"874\u00f75=\"], [\"741\u00f77=\", \"101\u00f75=\"], [\"532\u00f73=\", \"621\u00f77=\"], [\"803\u00f78=\", \"622\u00f75=\"], [\"945\u00f77=\", \"663\u00f77=\"], [\"356\u00f72=\", \"482\u00f74=\"], [\"930\u00f75=\", \"915\u00f72=\"], [\"127\u00f79=\", \"482\u00f73=\"], [\"746\u00f77=\", \"528\u00f76=\"], [\"848\u00f72=\", \"331\u00f79=\"], [\"736\u00f76=\", \"647\u00f78=\"], [\"712\u00f74=\", \"981\u00f74=\"], [\"655\u00f76=\", \"880\u00f75=\"], [\"594\u00f73=\", \"722\u00f73=\"], [\"351\u00f77=\", \"519\u00f78=\"], [\"519\u00f77=\", \"477\u00f77=\"], [\"782\u00f73=\", \"632\u00f75=\"], [\"114\u00f72=\", \"977\u00f78=\"], [\"537\u00f77=\", \"984\u00f74=\"], [\"819\u00f74=\", \"459\u00f76=\"], [\"231\u00f73=\", \"480\u00f79=\"], [\"873\u00f76=\", \"237\u00f73=\"]];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and every division problem's operands.\n# Each entry's Old text is unique across the document, so a single\n# Find/Execute per pair unambiguously targets the one run to change.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2024-08-04 Sunday\"; New = \"2024-08-05 Monday\" },\n    @{ Old = \"864\u00f76=\"; New = \"810\u00f77=\" },\n    @{ Old = \"252\u00f79=\"; New = \"447\u00f78=\" },\n    @{ Old = \"826\u00f79=\"; New = \"713\u00f77=\" },\n    @{ Old = \"390\u00f75=\"; New = \"874\u00f75=\" },\n    @{ Old = \"741\u00f77=\"; New = \"101\u00f75=\" },\n    @{ Old = \"532\u00f73=\"; New = \"621\u00f77=\" },\n    @{ Old = \"803\u00f78=\"; New = \"622\u00f75=\" },\n    @{ Old = \"945\u00f77=\"; New = \"663\u00f77=\" },\n    @{ Old = \"356\u00f72=\"; New = \"482\u00f74=\" },\n    @{ Old = \"930\u00f75=\"; New = \"915\u00f72=\" },\n    @{ Old = \"127\u00f79=\"; New = \"482\u00f73=\" },\n    @{ Old = \"746\u00f77=\"; New = \"528\u00f76=\" },\n    @{ Old = \"848\u00f72=\"; New = \"331\u00f79=\" },\n    @{ Old = \"736\u00f76=\"; New = \"647\u00f78=\" },\n    @{ Old = \"712\u00f74=\"; New = \"981\u00f74=\" },\n    @{ Old = \"655\u00f76=\"; New = \"880\u00f75=\" },\n    @{ Old = \"594\u00f73=\"; New = \"722\u00f73=\" },\n    @{ Old = \"351\u00f77=\"; New = \"519\u00f78=\" },\n    @{ Old = \"519\u00f77=\"; New = \"477\u00f77=\" },\n    @{ Old = \"782\u00f73=\"; New = \"632\u00f75=\" },\n    @{ Old = \"114\u00f72=\"; New = \"977\u00f78=\" },\n    @{ Old = \"537\u00f77=\"; New = \"984\u00f74=\" },\n    @{ Old = \"819\u00f74=\"; New = \"459\u00f76=\" },\n    @{ Old = \"231\u00f73=\"; New = \"480\u00f79=\" },\n    @{ Old = \"873\u00f76=\"; New = \"237\u00f73=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $ok = $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n    if (-not $ok) {\n        throw \"Find/Replace failed for: $($pair.Old)\"\n    }\n}\n"}
